# ----------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right before the "总计" (total) sheet and
#    populate it with that quarter's fund-holding table (same layout as the other
#    quarterly sheets: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值/仓位排名).
# 2) Insert a matching row at the top of the "总计" summary sheet's data table and
#    shift the previously-existing quarters down by one row.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Step 1: new "2022-Q1" worksheet, inserted immediately before "总计" ---------
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# Re-fetch "总计" by name now that the sheet collection has shifted: a handle
# captured before Add() keeps resolving by its old *position*, which after the
# insert now belongs to the brand-new sheet instead of "总计".
$totalSheet = $wb.Worksheets.Item("总计")

$headerStyle = $totalSheet.Range("B1").Style
$indexStyle  = $totalSheet.Range("A2").Style

$fundHeaders = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $fundHeaders[$col - 2]
    $newSheet.Cells.Item(1, $col).Style = $headerStyle
}

# Columns: 0:A-index 1:B-code 2:C-name 3:D-scale 4:E-stockPos 5:F-posPct 6:G-marketValue 7:H-rank
# B-G hold text (fund codes / decimal strings with significant trailing zeros) -> keep as text
$fundData = @(
    @(0, "010336", "中欧悦享生活混合A", "44.28", "90.44", "3.54", "1.5675", 9),
    @(1, "010787", "华安优势企业混合A", "36.73", "92.23", "3.99", "1.4655", 9),
    @(2, "004263", "华安沪港深机会灵活配置混合", "13.52", "94.79", "8.62", "1.1654", 1),
    @(3, "011162", "博时港股通领先趋势混合A", "23.56", "80.83", "3.47", "0.8175", 9),
    @(4, "010326", "博时消费创新混合A", "19.00", "82.83", "3.27", "0.6213", 10),
    @(5, "013414", "太平智远三个月定期开放股票", "8.69", "86.34", "4.08", "0.3546", 10),
    @(6, "009360", "招商创新增长混合A", "8.23", "91.96", "3.35", "0.2757", 10),
    @(7, "012987", "嘉合锦明混合A", "6.24", "63.80", "3.42", "0.2134", 6),
    @(8, "005521", "华安红利精选混合", "4.72", "91.71", "4.34", "0.2048", 8),
    @(9, "011238", "华安聚恒精选混合A", "5.14", "91.92", "3.98", "0.2046", 9),
    @(10, "010852", "中欧内需成长混合型证券投资基金A", "5.23", "91.46", "3.84", "0.2008", 9),
    @(11, "005620", "中欧品质消费股票A", "3.74", "90.47", "5.33", "0.1993", 9),
    @(12, "011163", "博时港股通领先趋势混合C", "4.68", "80.83", "3.47", "0.1624", 9),
    @(13, "012988", "嘉合锦明混合C", "3.82", "63.80", "3.42", "0.1306", 6),
    @(14, "010788", "华安优势企业混合C", "2.23", "92.23", "3.99", "0.0890", 9),
    @(15, "010327", "博时消费创新混合C", "2.60", "82.83", "3.27", "0.0850", 10),
    @(16, "006768", "华安沪港深优选混合", "0.84", "93.09", "8.78", "0.0738", 2),
    @(17, "005621", "中欧品质消费股票C", "1.11", "90.47", "5.33", "0.0592", 9),
    @(18, "010337", "中欧悦享生活混合C", "1.08", "90.44", "3.54", "0.0382", 9),
    @(19, "009361", "招商创新增长混合C", "1.04", "91.96", "3.35", "0.0348", 10),
    @(20, "011239", "华安聚恒精选混合C", "0.65", "91.92", "3.98", "0.0259", 9),
    @(21, "010853", "中欧内需成长混合型证券投资基金C", "0.67", "91.46", "3.84", "0.0257", 9),
)

for ($i = 0; $i -lt $fundData.Count; $i++) {
    $row = $fundData[$i]
    $r = $i + 2
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 1).Style = $indexStyle
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 3).Value = "'" + $row[2]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
}

$newSheet.Range("A1:H23").NumberFormat = "General"

# --- Step 2: push a new "2022-Q1" row onto the "总计" summary table --------------
$lastRow = 6   # existing data occupies rows 2..6 before this edit
for ($r = $lastRow; $r -ge 2; $r--) {
    $dest = $r + 1
    # column A is a 0-based running index ($dest row -> index $dest-2), not copied verbatim
    $totalSheet.Cells.Item($dest, 1).Value = $dest - 2
    $totalSheet.Cells.Item($dest, 1).Style = $totalSheet.Cells.Item($r, 1).Style
    $totalSheet.Cells.Item($dest, 2).Value = "'" + $totalSheet.Cells.Item($r, 2).Value2
    $totalSheet.Cells.Item($dest, 3).Value = $totalSheet.Cells.Item($r, 3).Value2
    $totalSheet.Cells.Item($dest, 4).Value = $totalSheet.Cells.Item($r, 4).Value2
}

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 1).Style = $indexStyle
$totalSheet.Cells.Item(2, 2).Value = "'2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 22
$totalSheet.Cells.Item(2, 4).Value = 8.02

$totalSheet.Range("A1:D7").NumberFormat = "General"
